$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4953.364
$ws.Range("I80").Value = 358.41177
$ws.Range("K80").Value = 1075.23531
$ws.Range("M80").Value = -77.23531000000003
$ws.Range("H83").Value = 4953.364
$ws.Range("I83").Value = 358.41177
$ws.Range("K83").Value = 3225.70593
$ws.Range("M83").Value = 1766.29407
$ws.Range("H92").Value = 62092012
$ws.Range("I92").Value = 3704146
$ws.Range("J92").Value = 500001000
$ws.Range("K92").Value = 3704146
$ws.Range("L92").Value = 500001000
$ws.Range("M92").Value = -3702898
$ws.Range("N92").Value = -500003496
$ws.Range("H100").Value = 27779526
$ws.Range("I100").Value = 33335032
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 33335032
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -33334491
$ws.Range("N100").Value = -3082
$ws.Range("H125").Value = 4381.6665
$ws.Range("I125").Value = 7921.3335
$ws.Range("J125").Value = 842
$ws.Range("K125").Value = 71292.0015
$ws.Range("L125").Value = 7578
$ws.Range("M125").Value = -68832.0015
$ws.Range("N125").Value = -12498
$ws.Range("H137").Value = 2270.2
$ws.Range("I137").Value = 1454.56
$ws.Range("J137").Value = 4309.3
$ws.Range("K137").Value = 4363.68
$ws.Range("L137").Value = 12927.9
$ws.Range("M137").Value = -1813.68
$ws.Range("N137").Value = -18027.9

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5657.143
$ws.Range("I3").Value = 1966.6666
$ws.Range("K3").Value = 1966.6666
$ws.Range("M3").Value = -1851.6666
$ws.Range("H32").Value = 9268.903
$ws.Range("I32").Value = 6808.8193
$ws.Range("J32").Value = 25371.273
$ws.Range("K32").Value = 6808.8193
$ws.Range("L32").Value = 25371.273
$ws.Range("M32").Value = -6521.8193
$ws.Range("N32").Value = -25945.273
$ws.Range("H45").Value = 6391.048
$ws.Range("I45").Value = 6978
$ws.Range("K45").Value = 6978
$ws.Range("M45").Value = -6601
$ws.Range("H61").Value = 208824.12
$ws.Range("I61").Value = 6094.769
$ws.Range("J61").Value = 437996.44
$ws.Range("K61").Value = 6094.769
$ws.Range("L61").Value = 437996.44
$ws.Range("M61").Value = -5882.769
$ws.Range("N61").Value = -438420.44
$ws.Range("H63").Value = 200021580
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 200021580
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1734.641
$ws.Range("I74").Value = 1360.1111
$ws.Range("J74").Value = 2577.3333
$ws.Range("K74").Value = 1360.1111
$ws.Range("L74").Value = 2577.3333
$ws.Range("M74").Value = -486.1111000000001
$ws.Range("N74").Value = -4325.3333
$ws.Range("H77").Value = 1734.641
$ws.Range("I77").Value = 1360.1111
$ws.Range("J77").Value = 2577.3333
$ws.Range("K77").Value = 6800.5555
$ws.Range("L77").Value = 12886.6665
$ws.Range("M77").Value = -2432.5555
$ws.Range("N77").Value = -21622.6665
$ws.Range("H102").Value = 2647139
$ws.Range("I102").Value = 3704934.5
$ws.Range("K102").Value = 3704934.5
$ws.Range("M102").Value = -3703312.5
$ws.Range("H132").Value = 3850348
$ws.Range("I132").Value = 2896.5789
$ws.Range("J132").Value = 14293430
$ws.Range("K132").Value = 8689.736699999999
$ws.Range("L132").Value = 42880290
$ws.Range("M132").Value = -6159.736699999999
$ws.Range("N132").Value = -42885350
$ws.Range("H136").Value = 208824.12
$ws.Range("I136").Value = 6094.769
$ws.Range("J136").Value = 437996.44
$ws.Range("K136").Value = 18284.307
$ws.Range("L136").Value = 1313989.32
$ws.Range("M136").Value = -15734.307
$ws.Range("N136").Value = -1319089.32

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = 4000
$ws.Range("K16").Value = 4000
$ws.Range("M16").Value = -3830
$ws.Range("H82").Value = 12618.286
$ws.Range("I82").Value = 3032
$ws.Range("J82").Value = 25400
$ws.Range("K82").Value = 3032
$ws.Range("L82").Value = 25400
$ws.Range("M82").Value = -2649
$ws.Range("N82").Value = -26166
$ws.Range("H85").Value = 12618.286
$ws.Range("I85").Value = 3032
$ws.Range("J85").Value = 25400
$ws.Range("K85").Value = 3032
$ws.Range("L85").Value = 25400
$ws.Range("M85").Value = -1706
$ws.Range("N85").Value = -28052
$ws.Range("H134").Value = 21381.105
$ws.Range("I134").Value = 4116.915
$ws.Range("J134").Value = 102522.8
$ws.Range("K134").Value = 12350.745
$ws.Range("L134").Value = 307568.4
$ws.Range("M134").Value = -9815.744999999999
$ws.Range("N134").Value = -312638.4

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 500
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H31").Value = 5231.2983
$ws.Range("I31").Value = 2154.611
$ws.Range("J31").Value = 6651.3076
$ws.Range("K31").Value = 2154.611
$ws.Range("L31").Value = 6651.3076
$ws.Range("M31").Value = -1859.611
$ws.Range("N31").Value = -7241.3076
$ws.Range("H34").Value = 5231.2983
$ws.Range("I34").Value = 2154.611
$ws.Range("J34").Value = 6651.3076
$ws.Range("K34").Value = 2154.611
$ws.Range("L34").Value = 6651.3076
$ws.Range("M34").Value = -1952.611
$ws.Range("N34").Value = -7055.3076
$ws.Range("H99").Value = 52333.332
$ws.Range("I99").Value = 52333.332
$ws.Range("K99").Value = 52333.332
$ws.Range("M99").Value = -50835.332
$ws.Range("H116").Value = 33750
$ws.Range("J116").Value = 33750
$ws.Range("L116").Value = 33750
$ws.Range("N116").Value = -42928
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -57258
$ws.Range("H126").Value = 52333.332
$ws.Range("I126").Value = 52333.332
$ws.Range("K126").Value = 156999.996
$ws.Range("M126").Value = -154529.996
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 2036.8966
$ws.Range("I132").Value = 1524.9131
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 4574.7393
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -2044.7393
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 280966.28
$ws.Range("I134").Value = 3313.5186
$ws.Range("J134").Value = 1113924.5
$ws.Range("K134").Value = 9940.5558
$ws.Range("L134").Value = 3341773.5
$ws.Range("M134").Value = -7405.5558
$ws.Range("N134").Value = -3346843.5

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 493.7
$ws.Range("I107").Value = 442.14285
$ws.Range("J107").Value = 614
$ws.Range("K107").Value = 1326.42855
$ws.Range("L107").Value = 1842
$ws.Range("M107").Value = 593.5714499999999
$ws.Range("N107").Value = -5682
$ws.Range("H113").Value = 1190994.4
$ws.Range("I113").Value = 1852329.2
$ws.Range("J113").Value = 591.8
$ws.Range("K113").Value = 5556987.6
$ws.Range("L113").Value = 1775.4
$ws.Range("M113").Value = -5554817.6
$ws.Range("N113").Value = -6115.4
$ws.Range("H132").Value = 1754.0714
$ws.Range("I132").Value = 760.8
$ws.Range("J132").Value = 2305.889
$ws.Range("K132").Value = 6847.2
$ws.Range("L132").Value = 20753.001
$ws.Range("M132").Value = -4317.2
$ws.Range("N132").Value = -25813.001

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3666.6667
$ws.Range("H80").Value = 5825.241
$ws.Range("I80").Value = 7594.1055
$ws.Range("J80").Value = 2464.4
$ws.Range("K80").Value = 7594.1055
$ws.Range("L80").Value = 2464.4
$ws.Range("M80").Value = -6596.1055
$ws.Range("N80").Value = -4460.4
$ws.Range("H83").Value = 5825.241
$ws.Range("I83").Value = 7594.1055
$ws.Range("J83").Value = 2464.4
$ws.Range("K83").Value = 37970.5275
$ws.Range("L83").Value = 12322
$ws.Range("M83").Value = -32978.5275
$ws.Range("N83").Value = -22306
$ws.Range("H102").Value = 2099.756
$ws.Range("I102").Value = 1850.6538
$ws.Range("J102").Value = 2531.5334
$ws.Range("K102").Value = 1850.6538
$ws.Range("L102").Value = 2531.5334
$ws.Range("M102").Value = -228.6538
$ws.Range("N102").Value = -5775.5334
$ws.Range("H126").Value = 9441.440000000001
$ws.Range("J126").Value = 2719
$ws.Range("L126").Value = 8157
$ws.Range("N126").Value = -13097
$ws.Range("H132").Value = 3128.7273
$ws.Range("I132").Value = 3001.5405
$ws.Range("J132").Value = 3291
$ws.Range("K132").Value = 9004.621500000001
$ws.Range("L132").Value = 9873
$ws.Range("M132").Value = -6474.621500000001
$ws.Range("N132").Value = -14933

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 655700.3
$ws.Range("I82").Value = 1667650.4
$ws.Range("J82").Value = 103727.55
$ws.Range("K82").Value = 1667650.4
$ws.Range("L82").Value = 103727.55
$ws.Range("M82").Value = -1667289.4
$ws.Range("N82").Value = -104449.55
$ws.Range("H85").Value = 655700.3
$ws.Range("I85").Value = 1667650.4
$ws.Range("J85").Value = 103727.55
$ws.Range("K85").Value = 1667650.4
$ws.Range("L85").Value = 103727.55
$ws.Range("M85").Value = -1666402.4
$ws.Range("N85").Value = -106223.55

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 939.3
$ws.Range("I126").Value = 821.44446
$ws.Range("K126").Value = 2464.33338
$ws.Range("M126").Value = 5.666619999999966
$ws.Range("H132").Value = 1303.8868
$ws.Range("I132").Value = 1013
$ws.Range("J132").Value = 3215.4285
$ws.Range("K132").Value = 3039
$ws.Range("L132").Value = 9646.2855
$ws.Range("M132").Value = -509
$ws.Range("N132").Value = -14706.2855
$ws.Range("H136").Value = 2722.5
$ws.Range("I136").Value = 2951.923
$ws.Range("J136").Value = 2509.4644
$ws.Range("K136").Value = 8855.769
$ws.Range("L136").Value = 7528.3932
$ws.Range("M136").Value = -6305.769
$ws.Range("N136").Value = -12628.3932
